$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.029.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4326"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3711"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07396"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9309"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.008.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.716"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06871"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.017.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.104"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.197.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.048"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.462"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.686"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08981"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8053"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.748"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.171"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05492"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.119"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.007"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5226"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.982"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1687"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.725"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06709"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4863"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.868"
$ws.Range("D51").Style = "Normal"

$ws.Range("E2").Value = "  -3.50%  "
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -5.87%  "
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("E10").Value = "  -4.99%  "
$ws.Range("E11").Value = "  -6.53%  "
$ws.Range("E12").Value = "  +5.35%  "
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -4.86%  "
$ws.Range("E18").Value = "  -5.65%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -5.88%  "
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("E22").Value = "  -4.38%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("E30").Value = "  -8.01%  "
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("E32").Value = "  -6.12%  "
$ws.Range("E33").Value = "  -7.04%  "
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("E40").Value = "  -5.02%  "
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("E42").Value = "  -6.70%  "
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("E44").Value = "  -6.75%  "
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("E46").Value = "  -6.19%  "
$ws.Range("E47").Value = "  -6.96%  "
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  -5.47%  "
$ws.Range("E51").Value = "  -14.88%  "
